$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.723.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.697.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  +0.35%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.17%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4079'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  -2.50%  '
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '52.88'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08919'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.266'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.029'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001323'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.698.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07045'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.004'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.98%  '
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.94%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.699.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").Value = '  +8.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.366'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '136.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.164'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.539'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08709'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.057'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.096'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.83%  '
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '14.50'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.884'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09239'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.475'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7669'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.592'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7182'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.225'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '140.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.320'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07989'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.14%  '
